# CYRS_Review.xlsx cross-review update
# - Reviewer fills in the "Decision" (G) and "Comment" (I) columns for the
#   seven open points logged on the "Cross review points " sheet.
# - View state (selection / scroll / zoom) and a couple of cosmetic
#   row-height / column-width tweaks follow, mirroring the author's session.

$wb = $excel.ActiveWorkbook

$wsIntro  = $wb.Worksheets.Item("Introduction ")
$wsReview = $wb.Worksheets.Item("Cross review points ")

# ---------------------------------------------------------------------
# Cross review points: Decision (G) + Comment (I) for rows 2-8
# ---------------------------------------------------------------------

# Row 2 - Document status table already covered by the document history
$wsReview.Range("G2").Value = "Accepted"
$wsReview.Range("I2").Value = "This is already included in the document history"

# Row 3 - Document should stay Draft + history moved to top of the document
$wsReview.Range("G3").Value = "Accepted"
$wsReview.Range("I3").Value = "Modified, and document history is moved to the start of the document"

# Row 4 - Missing reference table for the reference document
$wsReview.Range("G4").Value = "Accepted"
$wsReview.Range("I4").Value = "The reference table that references each requirement to its parent requirement in the CRS document is at the end of the CYRS document"

# Row 5 - CYRS requirements should read at system level, not SW level
$wsReview.Range("G5").Value = "Accepted"
$wsReview.Range("I5").Value = "Modified as a system and not software requirements"

# Row 6 - HW requirement detail belongs to the HSI document
$wsReview.Range("G6").Value = "Refused"
$wsReview.Range("I6").Value = "As the HSI is the document related to the chosen HW components so in the CYRS for example a displayer in the high level and in the HSI its specified to a 16*2 lcd "

# Row 7 - Buzzer HW/SW relation is expected to live across HSI/SRS
$wsReview.Range("G7").Value = "Refused"
$wsReview.Range("I7").Value = 'Thiss requirement describes a high level relation between the HW component""Buzzer"  that will be referenced to the HIS and a SW requirement"That buzzer is on when specific time is hit" that will be referenced to the SRS '

# Row 8 - All requirements updated to system level perspective
$wsReview.Range("G8").Value = "Accepted"
$wsReview.Range("I8").Value = "Modified"

# ---------------------------------------------------------------------
# Row heights on "Cross review points " now that the Comment column has
# wrapped, multi-line text in it (auto-fit in the authoring session).
# ---------------------------------------------------------------------
$wsReview.Rows.Item(3).RowHeight = 75
$wsReview.Rows.Item(4).RowHeight = 60
$wsReview.Rows.Item(5).RowHeight = 90
$wsReview.Rows.Item(6).RowHeight = 75
$wsReview.Rows.Item(7).RowHeight = 105
$wsReview.Rows.Item(8).RowHeight = 30

# Comment column widened to fit the new text
$wsReview.Columns.Item(9).ColumnWidth = 39

# ---------------------------------------------------------------------
# View state: scroll / zoom / selection, matching the end of the session
# ---------------------------------------------------------------------
$wsReview.Activate()
$excel.ActiveWindow.Zoom = 70
$wsReview.Range("L5").Select()
$wsReview.Range("E1").Select()

$wsIntro.Activate()
$wsIntro.Range("D9:H9").Select()

$wsReview.Activate()
